# New .ttl from Google sheet has been generated
# Column D ("rdf:type(separator=",")") for the iop:VariableSet / skos:Concept
# rows no longer also declares skos:Concept - these rows are now plain
# iop:VariableSet entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,94,95,96,97,98,99,100,101,102,103,104,105,106)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "iop:VariableSet"
}
